$wb = $excel.ActiveWorkbook

# --- 1. Add "test" entry to the struct_list sheet ---
$structList = $wb.Worksheets.Item("struct_list")
$structList.Range("A7").Value = "test"

# --- 2. Create the new "test" worksheet by copying the "template" sheet,
#        placed after the last sheet ("service") ---
$template = $wb.Worksheets.Item("template")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)

$testSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$testSheet.Name = "test"

# --- 3. Fill in the register rows for the "test" struct (array variants).
#        Column A is filled top-to-bottom first (skipping row 12, the
#        "double" row, which is appended later), then column C, then the
#        "double" row's name cell last -- this reproduces the exact order
#        the new shared strings were authored in.
$testSheet.Range("A3").Value = "arr_u8"
$testSheet.Range("A4").Value = "arr_u16"
$testSheet.Range("A5").Value = "arr_u32"
$testSheet.Range("A6").Value = "arr_u64"
$testSheet.Range("A7").Value = "arr_s8"
$testSheet.Range("A8").Value = "arr_s16"
$testSheet.Range("A9").Value = "arr_s32"
$testSheet.Range("A10").Value = "arr_s64"
$testSheet.Range("A11").Value = "arr_float"
$testSheet.Range("A13").Value = "arr_char"

$testSheet.Range("C3").Value = "u8"
$testSheet.Range("C4").Value = "u16"
$testSheet.Range("C5").Value = "u32"
$testSheet.Range("C6").Value = "u64"
$testSheet.Range("C7").Value = "s8"
$testSheet.Range("C8").Value = "s16"
$testSheet.Range("C9").Value = "s32"
$testSheet.Range("C10").Value = "s64"
$testSheet.Range("C11").Value = "float"
$testSheet.Range("C12").Value = "double"
$testSheet.Range("C13").Value = "char"

$testSheet.Range("A12").Value = "arr_double"

# --- 4. array_len / p_def columns ---
$testSheet.Range("D3").Value = 20
$testSheet.Range("D4").Value = 20
$testSheet.Range("D5").Value = 20
$testSheet.Range("D6").Value = 20
$testSheet.Range("D7").Value = 20
$testSheet.Range("D8").Value = 20
$testSheet.Range("D9").Value = 20
$testSheet.Range("D10").Value = 20
$testSheet.Range("D11").Value = 20
$testSheet.Range("D12").Value = 20
$testSheet.Range("D13").Value = 20

$testSheet.Range("F3").Value = 110
$testSheet.Range("F4").Value = "auto"
$testSheet.Range("F5").Value = "auto"
$testSheet.Range("F6").Value = "auto"
$testSheet.Range("F7").Value = "auto"
$testSheet.Range("F8").Value = "auto"
$testSheet.Range("F9").Value = "auto"
$testSheet.Range("F10").Value = "auto"
$testSheet.Range("F11").Value = "auto"
$testSheet.Range("F12").Value = "auto"
$testSheet.Range("F13").Value = "auto"

# --- 5. Restore view/selection state to match the final snapshot ---
$testSheet.Range("F4").Select() | Out-Null

$structList.Range("A7").Select() | Out-Null

$osSheet = $wb.Worksheets.Item("os")
$osSheet.Range("Q14").Select() | Out-Null

$testSheet.Activate() | Out-Null
